$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 13.31988742298459
    "C2" = 7.450744628479877
    "D2" = 8.521319697582616
    "F2" = 37.08832220141597
    "G2" = 3.688117323192412
    "I2" = 27.92132695538892
    "J2" = 10.53945038367646
    "K2" = 10.4280967869138
    "L2" = 11.67411010394762
    "O2" = 28.39819587319363
    "B3" = 13.08552156146055
    "C3" = 7.411935438787108
    "D3" = 8.499367562941936
    "F3" = 37.19029826516469
    "G3" = 3.689904212881955
    "I3" = 28.02635645040964
    "J3" = 10.56173771461897
    "K3" = 10.26033729571833
    "L3" = 11.66928250930501
    "O3" = 28.49563815991452
    "B4" = 12.94137273127791
    "C4" = 7.387791827781644
    "D4" = 8.487012254949335
    "F4" = 37.26023124280371
    "G4" = 3.691060133755287
    "I4" = 28.09540865967248
    "J4" = 10.57630978554013
    "K4" = 10.15723738312069
    "L4" = 11.66781652829422
    "O4" = 28.56053386295494
    "B5" = 12.8826429082736
    "C5" = 7.377877156282215
    "D5" = 8.482263365854205
    "F5" = 37.29056766527374
    "G5" = 3.691546002718167
    "I5" = 28.12469579483335
    "J5" = 10.58247169344236
    "K5" = 10.11524896409019
    "L5" = 11.66759744746888
    "O5" = 28.58825236498105
    "B6" = 12.87289378060865
    "C6" = 7.376226337238126
    "D6" = 8.481492188459884
    "F6" = 37.29571595439357
    "G6" = 3.69162757742421
    "I6" = 28.12962823270574
    "J6" = 10.58350839705342
    "K6" = 10.10827982444378
    "L6" = 11.66758396254298
    "O6" = 28.59293186077539
    "B7" = 12.94058053292107
    "C7" = 7.387658418425537
    "D7" = 8.486947047283957
    "F7" = 37.26063293095336
    "G7" = 3.691066626280639
    "I7" = 28.0957989885131
    "J7" = 10.57639198090186
    "K7" = 10.15667094303118
    "L7" = 11.66781203989444
    "O7" = 28.5609025312727
    "B8" = 13.23917348349084
    "C8" = 7.437430410185105
    "D8" = 8.513519612208361
    "F8" = 37.12196309330306
    "G8" = 3.688721273726432
    "I8" = 27.95659406979993
    "J8" = 10.54695113517174
    "K8" = 10.37030251111157
    "L8" = 11.67213559165877
    "O8" = 28.43074207671577
    "B9" = 13.8195683975515
    "C9" = 7.532399812801029
    "D9" = 8.574379071483897
    "F9" = 36.90820481780798
    "G9" = 3.684586238115892
    "I9" = 27.71981394975688
    "J9" = 10.49623818382338
    "K9" = 10.78634421517178
    "L9" = 11.69242740718602
    "O9" = 28.21572943381213
    "B10" = 14.23872629625061
    "C10" = 7.600387531766915
    "D10" = 8.624196895970773
    "F10" = 36.78674485202183
    "G10" = 3.681828317484009
    "I10" = 27.56790560003646
    "J10" = 10.46322893584884
    "K10" = 11.08746481898427
    "L10" = 11.71443188785597
    "O10" = 28.0823357786628
    "B11" = 14.42701346220322
    "C11" = 7.630890682021235
    "D11" = 8.647915494312043
    "F11" = 36.73923971203355
    "G11" = 3.680633871476543
    "I11" = 27.5035848467654
    "J11" = 10.44912852982537
    "K11" = 11.2229041778308
    "L11" = 11.72595784891356
    "O11" = 28.02699848531423
    "B12" = 14.49790470114239
    "C12" = 7.642377051162069
    "D12" = 8.657044174733555
    "F12" = 36.7223663134578
    "G12" = 3.68019016823187
    "I12" = 27.47991603084546
    "J12" = 10.44392025949095
    "K12" = 11.27392534783486
    "L12" = 11.73053802490182
    "O12" = 28.00681308197582
    "B13" = 14.48265620072335
    "C13" = 7.639906179583372
    "D13" = 8.655071694115488
    "F13" = 36.72595065384564
    "G13" = 3.680285345496674
    "I13" = 27.48498292802636
    "J13" = 10.44503612314296
    "K13" = 11.2629495864581
    "L13" = 11.72954205626008
    "O13" = 28.01112612234613
    "B14" = 14.43285428696038
    "C14" = 7.631836984995105
    "D14" = 8.648663599816567
    "F14" = 36.73782915708953
    "G14" = 3.680597195470984
    "I14" = 27.50162380780777
    "J14" = 10.44869741431161
    "K14" = 11.22710731756024
    "L14" = 11.72633035748277
    "O14" = 28.02532238905999
    "B15" = 14.40229402727026
    "C15" = 7.626885867604773
    "D15" = 8.64475744676834
    "F15" = 36.74525043430928
    "G15" = 3.680789332318304
    "I15" = 27.5119064404565
    "J15" = 10.4509571406385
    "K15" = 11.20511684858862
    "L15" = 11.72439109421551
    "O15" = 28.03411827425775
    "B16" = 14.22636753306823
    "C16" = 7.598385237808874
    "D16" = 8.622667667017041
    "F16" = 36.79000547491948
    "G16" = 3.681907584084619
    "I16" = 27.57220538192059
    "J16" = 10.46416881972256
    "K16" = 11.07857861401402
    "L16" = 11.71370893853598
    "O16" = 28.08605983750026
    "B17" = 14.1177826640253
    "C17" = 7.580789779654117
    "D17" = 8.609383381381429
    "F17" = 36.81944685989586
    "G17" = 3.6826089705861
    "I17" = 27.61042207055696
    "J17" = 10.47250797776632
    "K17" = 11.00052389605249
    "L17" = 11.70754236690765
    "O17" = 28.11929389333329
    "B18" = 14.05510627777051
    "C18" = 7.570629556964685
    "D18" = 8.601842334225289
    "F18" = 36.83710992324401
    "G18" = 3.683018053220103
    "I18" = 27.63285349192425
    "J18" = 10.47739065286377
    "K18" = 10.95548632610593
    "L18" = 11.70413836830444
    "O18" = 28.13891219291159
    "B19" = 14.03384917812091
    "C19" = 7.56718275243283
    "D19" = 8.599306339698598
    "F19" = 36.84321551145264
    "G19" = 3.683157535579129
    "I19" = 27.64052569963828
    "J19" = 10.47905866274248
    "K19" = 10.94021429943906
    "L19" = 11.70301044004404
    "O19" = 28.14564095879522
    "B20" = 14.12936509049134
    "C20" = 7.582666988220321
    "D20" = 8.610787230444865
    "F20" = 36.81623729623015
    "G20" = 3.682533720927896
    "I20" = 27.60630724379456
    "J20" = 10.47161134063084
    "K20" = 11.00884804127234
    "L20" = 11.70818404185863
    "O20" = 28.11570400954724
    "B21" = 14.44749392124531
    "C21" = 7.634208881089727
    "D21" = 8.650541864027247
    "F21" = 36.73430985911214
    "G21" = 3.680505364296606
    "I21" = 27.49671730222256
    "J21" = 10.44761844568763
    "K21" = 11.23764264988379
    "L21" = 11.72726788172151
    "O21" = 28.02113170490254
    "B22" = 14.65299748793857
    "C22" = 7.667516887937073
    "D22" = 8.677378018614498
    "F22" = 36.6872696785166
    "G22" = 3.67922986833493
    "I22" = 27.4291046005321
    "J22" = 10.43270251837725
    "K22" = 11.38559874937156
    "L22" = 11.74099539484928
    "O22" = 27.963809733511
    "B23" = 14.54355778180208
    "C23" = 7.649775455001759
    "D23" = 8.662978568552578
    "F23" = 36.71178029418643
    "G23" = 3.67990604917033
    "I23" = 27.46482370642772
    "J23" = 10.44059358780327
    "K23" = 11.30679020231286
    "L23" = 11.73355476699355
    "O23" = 27.99399266220895
    "B24" = 14.12412944761197
    "C24" = 7.58181844014873
    "D24" = 8.610152250068214
    "F24" = 36.81768604494694
    "G24" = 3.682567723084704
    "I24" = 27.60816612367969
    "J24" = 10.47201643494632
    "K24" = 11.00508519908183
    "L24" = 11.70789350046994
    "O24" = 28.11732540248512
    "B25" = 13.66354191011074
    "C25" = 7.507008659012501
    "D25" = 8.557001095495401
    "F25" = 36.95979057798991
    "G25" = 3.685655480220668
    "I25" = 27.77999537909469
    "J25" = 10.50920894124551
    "K25" = 10.67439590366976
    "L25" = 11.69242740718602
    "O25" = 28.26958376040688
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
